# "Enabled all the tests."
#
# On the "Test Cases" sheet, column D ("Runmode") holds "Y"/"N" flags for
# whether each test case is run. Every row (D2:D89) except D81 -- which was
# already "Y" -- is currently "N". Flip them all to "Y" so every test runs,
# and leave the already-enabled D81 alone. The now-unreferenced "N" shared
# string is dropped automatically when the workbook is saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

for ($r = 2; $r -le 89; $r++) {
    if ($r -eq 81) { continue }
    $ws.Cells.Item($r, 4).Value = "Y"
}

# Mirror the author's final selection: the whole Runmode column they just
# edited, anchored back at the top.
$ws.Range("D2:D89").Select()
